$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Stable accommodation.global.sdf"
$ws.Range("C1").Value = "Unstable accommodation and/or homeless.global.sdf"
$ws.Range("D1").Value = "In detention.global.sdf"
$ws.Range("E1").Value = "Other.global.sdf"
$ws.Range("F1").Value = "Not known / missing.global.sdf"
$ws.Range("G1").Value = "Total.global.sdf"
